$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 436, shifting existing rows 436-492 down to 437-493
$ws.Rows.Item(436).Insert()

# Populate the new row 436 with data
$ws.Cells.Item(436, 1).Value = 4
$ws.Cells.Item(436, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(436, 3).Value = "Los Lagos"
$ws.Cells.Item(436, 4).Value = 45131
$ws.Cells.Item(436, 5).Value = 10
$ws.Cells.Item(436, 6).Value = 100112003
$ws.Cells.Item(436, 7).Value = "Ajo"
$ws.Cells.Item(436, 8).Value = "Chino"
$ws.Cells.Item(436, 9).Value = "Primera"
$ws.Cells.Item(436, 10).Value = 120
$ws.Cells.Item(436, 11).Value = 22000
$ws.Cells.Item(436, 12).Value = 22000
$ws.Cells.Item(436, 13).Value = 22000
$ws.Cells.Item(436, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(436, 15).Value = "China"
$ws.Cells.Item(436, 16).Value = 2200
$ws.Cells.Item(436, 17).Value = 10
$ws.Cells.Item(436, 18).Value = "Hortaliza"
